# Fixed some minor new clean install bugs
# - DeviceDiscovery (B2) toggled off
# - DeviceDiscoverySSH (B3) toggled on
# - Selection left on C15 (scrolled back to top of sheet)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Config")

# Toggle the two boolean config values
$ws.Range("B2").Value = $false   # DeviceDiscovery      -> False
$ws.Range("B3").Value = $true    # DeviceDiscoverySSH    -> True

# Restore the view: top of sheet visible, C15 selected/active
$ws.Activate() | Out-Null
$ws.Range("A1").Select() | Out-Null
$ws.Range("C15").Select() | Out-Null
